# Apply edit: insert 3 new data rows before existing row 170 (Femacal de La
# Calera - Chirimoya weekly price data), shifting all subsequent rows down
# by 3 (new dimension becomes A1:T241).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three blank rows at row 170; existing rows 170.. shift down to 173..
$ws.Range("A170:T172").EntireRow.Insert()

# Common (unchanged across the three new rows) field values
$mercadoId   = 3
$mercado     = "Femacal de La Calera"
$region      = "Coquimbo"
$fecha       = 44837
$codreg      = 5
$tipo        = "Fruta"
$productoId  = 100107
$producto    = "Otros"
$categoriaId = 100107002
$categoria   = "Chirimoya"
$variedad    = "Cultivar IV Región"
$unidad      = "$/bandeja 10 kilos"
$origen      = "Provincia del Elquí"
$kgUnidad    = 10

# Row 170: Calidad "Especial"
$r = 170
$ws.Cells.Item($r, 1).Value  = $mercadoId
$ws.Cells.Item($r, 2).Value  = $mercado
$ws.Cells.Item($r, 3).Value  = $region
$ws.Cells.Item($r, 4).Value  = $fecha
$ws.Cells.Item($r, 5).Value  = $codreg
$ws.Cells.Item($r, 6).Value  = $tipo
$ws.Cells.Item($r, 7).Value  = $productoId
$ws.Cells.Item($r, 8).Value  = $producto
$ws.Cells.Item($r, 9).Value  = $categoriaId
$ws.Cells.Item($r, 10).Value = $categoria
$ws.Cells.Item($r, 11).Value = $variedad
$ws.Cells.Item($r, 12).Value = "Especial"
$ws.Cells.Item($r, 13).Value = 56
$ws.Cells.Item($r, 14).Value = 30000
$ws.Cells.Item($r, 15).Value = 30000
$ws.Cells.Item($r, 16).Value = 30000
$ws.Cells.Item($r, 17).Value = $unidad
$ws.Cells.Item($r, 18).Value = $origen
$ws.Cells.Item($r, 19).Value = 3000
$ws.Cells.Item($r, 20).Value = $kgUnidad

# Row 171: Calidad "Primera"
$r = 171
$ws.Cells.Item($r, 1).Value  = $mercadoId
$ws.Cells.Item($r, 2).Value  = $mercado
$ws.Cells.Item($r, 3).Value  = $region
$ws.Cells.Item($r, 4).Value  = $fecha
$ws.Cells.Item($r, 5).Value  = $codreg
$ws.Cells.Item($r, 6).Value  = $tipo
$ws.Cells.Item($r, 7).Value  = $productoId
$ws.Cells.Item($r, 8).Value  = $producto
$ws.Cells.Item($r, 9).Value  = $categoriaId
$ws.Cells.Item($r, 10).Value = $categoria
$ws.Cells.Item($r, 11).Value = $variedad
$ws.Cells.Item($r, 12).Value = "Primera"
$ws.Cells.Item($r, 13).Value = 57
$ws.Cells.Item($r, 14).Value = 27000
$ws.Cells.Item($r, 15).Value = 27000
$ws.Cells.Item($r, 16).Value = 27000
$ws.Cells.Item($r, 17).Value = $unidad
$ws.Cells.Item($r, 18).Value = $origen
$ws.Cells.Item($r, 19).Value = 2700
$ws.Cells.Item($r, 20).Value = $kgUnidad

# Row 172: Calidad "Segunda"
$r = 172
$ws.Cells.Item($r, 1).Value  = $mercadoId
$ws.Cells.Item($r, 2).Value  = $mercado
$ws.Cells.Item($r, 3).Value  = $region
$ws.Cells.Item($r, 4).Value  = $fecha
$ws.Cells.Item($r, 5).Value  = $codreg
$ws.Cells.Item($r, 6).Value  = $tipo
$ws.Cells.Item($r, 7).Value  = $productoId
$ws.Cells.Item($r, 8).Value  = $producto
$ws.Cells.Item($r, 9).Value  = $categoriaId
$ws.Cells.Item($r, 10).Value = $categoria
$ws.Cells.Item($r, 11).Value = $variedad
$ws.Cells.Item($r, 12).Value = "Segunda"
$ws.Cells.Item($r, 13).Value = 48
$ws.Cells.Item($r, 14).Value = 24000
$ws.Cells.Item($r, 15).Value = 24000
$ws.Cells.Item($r, 16).Value = 24000
$ws.Cells.Item($r, 17).Value = $unidad
$ws.Cells.Item($r, 18).Value = $origen
$ws.Cells.Item($r, 19).Value = 2400
$ws.Cells.Item($r, 20).Value = $kgUnidad
